# Fruta / hortaliza, semanal
# Insert a new weekly record at row 16 (Papaya, Vega Modelo de Temuco),
# pushing the existing data rows (old 16..78) down by one row to (17..79).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 16.
$ws.Rows("16:16").Insert()

# Populate the newly inserted row 16 with the new data point.
$ws.Cells.Item(16,1).Value  = 10
$ws.Cells.Item(16,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(16,3).Value  = "La Araucanía"
$ws.Cells.Item(16,4).Value  = 44676
$ws.Cells.Item(16,5).Value  = 9
$ws.Cells.Item(16,6).Value  = "Fruta"
$ws.Cells.Item(16,7).Value  = 100108
$ws.Cells.Item(16,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(16,9).Value  = 100108004
$ws.Cells.Item(16,10).Value = "Papaya"
$ws.Cells.Item(16,11).Value = "Sin especificar"
$ws.Cells.Item(16,12).Value = "Primera"
$ws.Cells.Item(16,13).Value = 30
$ws.Cells.Item(16,14).Value = 2500
$ws.Cells.Item(16,15).Value = 2500
$ws.Cells.Item(16,16).Value = 2500
$ws.Cells.Item(16,17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(16,18).Value = "Provincia del Elquí"
$ws.Cells.Item(16,19).Value = 2500
$ws.Cells.Item(16,20).Value = 1
